$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("2025-02-12 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-13 Thursday", 2)
$null = $d.Content.Find.Execute("35-1=34", $true, $false, $false, $false, $false, $true, 1, $false, "8-3=5", 2)
$null = $d.Content.Find.Execute("29-10=19", $true, $false, $false, $false, $false, $true, 1, $false, "57+42=99", 2)
$null = $d.Content.Find.Execute("77-5=72", $true, $false, $false, $false, $false, $true, 1, $false, "16-11=5", 2)
$null = $d.Content.Find.Execute("18-9=9", $true, $false, $false, $false, $false, $true, 1, $false, "10+13=23", 2)
$null = $d.Content.Find.Execute("46+17=63", $true, $false, $false, $false, $false, $true, 1, $false, "88-17=71", 2)
$null = $d.Content.Find.Execute("87-50=37", $true, $false, $false, $false, $false, $true, 1, $false, "9+16=25", 2)
$null = $d.Content.Find.Execute("27+26=53", $true, $false, $false, $false, $false, $true, 1, $false, "12+10=22", 2)
$null = $d.Content.Find.Execute("0+56=56", $true, $false, $false, $false, $false, $true, 1, $false, "79-61=18", 2)
$null = $d.Content.Find.Execute("61-41=20", $true, $false, $false, $false, $false, $true, 1, $false, "71-32=39", 2)
$null = $d.Content.Find.Execute("15+62=77", $true, $false, $false, $false, $false, $true, 1, $false, "36+16=52", 2)
$null = $d.Content.Find.Execute("12-0=12", $true, $false, $false, $false, $false, $true, 1, $false, "63-21=42", 2)
$null = $d.Content.Find.Execute("46+20=66", $true, $false, $false, $false, $false, $true, 1, $false, "57+34=91", 2)
$null = $d.Content.Find.Execute("36+51=87", $true, $false, $false, $false, $false, $true, 1, $false, "63-48=15", 2)
$null = $d.Content.Find.Execute("96-74=22", $true, $false, $false, $false, $false, $true, 1, $false, "35-35=0", 2)
$null = $d.Content.Find.Execute("12+2=14", $true, $false, $false, $false, $false, $true, 1, $false, "11-2=9", 2)
$null = $d.Content.Find.Execute("61-30=31", $true, $false, $false, $false, $false, $true, 1, $false, "99-52=47", 2)
$null = $d.Content.Find.Execute("44-32=12", $true, $false, $false, $false, $false, $true, 1, $false, "1+13=14", 2)
$null = $d.Content.Find.Execute("24-14=10", $true, $false, $false, $false, $false, $true, 1, $false, "13+33=46", 2)
$null = $d.Content.Find.Execute("26+65=91", $true, $false, $false, $false, $false, $true, 1, $false, "40-11=29", 2)
$null = $d.Content.Find.Execute("40-24=16", $true, $false, $false, $false, $false, $true, 1, $false, "45+6=51", 2)
$null = $d.Content.Find.Execute("88-72=16", $true, $false, $false, $false, $false, $true, 1, $false, "50+19=69", 2)
$null = $d.Content.Find.Execute("78-62=16", $true, $false, $false, $false, $false, $true, 1, $false, "38+33=71", 2)
$null = $d.Content.Find.Execute("40+59=99", $true, $false, $false, $false, $false, $true, 1, $false, "52-46=6", 2)
$null = $d.Content.Find.Execute("61-2=59", $true, $false, $false, $false, $false, $true, 1, $false, "23+55=78", 2)
$null = $d.Content.Find.Execute("46+5=51", $true, $false, $false, $false, $false, $true, 1, $false, "48+35=83", 2)
$null = $d.Content.Find.Execute("36+40=76", $true, $false, $false, $false, $false, $true, 1, $false, "54+22=76", 2)
$null = $d.Content.Find.Execute("83+3=86", $true, $false, $false, $false, $false, $true, 1, $false, "34-4=30", 2)
$null = $d.Content.Find.Execute("12+33=45", $true, $false, $false, $false, $false, $true, 1, $false, "24+44=68", 2)
$null = $d.Content.Find.Execute("57-41=16", $true, $false, $false, $false, $false, $true, 1, $false, "52-43=9", 2)
$null = $d.Content.Find.Execute("57+37=94", $true, $false, $false, $false, $false, $true, 1, $false, "95-68=27", 2)
$null = $d.Content.Find.Execute("23+72=95", $true, $false, $false, $false, $false, $true, 1, $false, "54+13=67", 2)
$null = $d.Content.Find.Execute("29-6=23", $true, $false, $false, $false, $false, $true, 1, $false, "70-62=8", 2)
$null = $d.Content.Find.Execute("95-15=80", $true, $false, $false, $false, $false, $true, 1, $false, "5+74=79", 2)
$null = $d.Content.Find.Execute("58+19=77", $true, $false, $false, $false, $false, $true, 1, $false, "86+10=96", 2)
$null = $d.Content.Find.Execute("33+52=85", $true, $false, $false, $false, $false, $true, 1, $false, "55-12=43", 2)
$null = $d.Content.Find.Execute("92-69=23", $true, $false, $false, $false, $false, $true, 1, $false, "47-45=2", 2)
$null = $d.Content.Find.Execute("59-43=16", $true, $false, $false, $false, $false, $true, 1, $false, "91-13=78", 2)
$null = $d.Content.Find.Execute("38+60=98", $true, $false, $false, $false, $false, $true, 1, $false, "41-4=37", 2)
$null = $d.Content.Find.Execute("95-77=18", $true, $false, $false, $false, $false, $true, 1, $false, "71-54=17", 2)
$null = $d.Content.Find.Execute("37+48=85", $true, $false, $false, $false, $false, $true, 1, $false, "61-50=11", 2)
$null = $d.Content.Find.Execute("22+52=74", $true, $false, $false, $false, $false, $true, 1, $false, "9+73=82", 2)
$null = $d.Content.Find.Execute("15+76=91", $true, $false, $false, $false, $false, $true, 1, $false, "88-59=29", 2)
$null = $d.Content.Find.Execute("71-29=42", $true, $false, $false, $false, $false, $true, 1, $false, "60-32=28", 2)
$null = $d.Content.Find.Execute("20+75=95", $true, $false, $false, $false, $false, $true, 1, $false, "48+1=49", 2)
$null = $d.Content.Find.Execute("4+19=23", $true, $false, $false, $false, $false, $true, 1, $false, "4+74=78", 2)
$null = $d.Content.Find.Execute("45+34=79", $true, $false, $false, $false, $false, $true, 1, $false, "77+22=99", 2)
$null = $d.Content.Find.Execute("30+64=94", $true, $false, $false, $false, $false, $true, 1, $false, "91-77=14", 2)
$null = $d.Content.Find.Execute("27+14=41", $true, $false, $false, $false, $false, $true, 1, $false, "31-9=22", 2)
$null = $d.Content.Find.Execute("63+22=85", $true, $false, $false, $false, $false, $true, 1, $false, "26+12=38", 2)
$null = $d.Content.Find.Execute("7+74=81", $true, $false, $false, $false, $false, $true, 1, $false, "97-93=4", 2)
$null = $d.Content.Find.Execute("63-9=54", $true, $false, $false, $false, $false, $true, 1, $false, "67+26=93", 2)
$null = $d.Content.Find.Execute("49+26=75", $true, $false, $false, $false, $false, $true, 1, $false, "74-61=13", 2)
$null = $d.Content.Find.Execute("92-39=53", $true, $false, $false, $false, $false, $true, 1, $false, "70-64=6", 2)
$null = $d.Content.Find.Execute("3-2=1", $true, $false, $false, $false, $false, $true, 1, $false, "92+7=99", 2)
$null = $d.Content.Find.Execute("6+11=17", $true, $false, $false, $false, $false, $true, 1, $false, "58-11=47", 2)
$null = $d.Content.Find.Execute("42+54=96", $true, $false, $false, $false, $false, $true, 1, $false, "61-47=14", 2)
$null = $d.Content.Find.Execute("41+40=81", $true, $false, $false, $false, $false, $true, 1, $false, "76-55=21", 2)
$null = $d.Content.Find.Execute("55-48=7", $true, $false, $false, $false, $false, $true, 1, $false, "27+7=34", 2)
$null = $d.Content.Find.Execute("42+57=99", $true, $false, $false, $false, $false, $true, 1, $false, "67+18=85", 2)
$null = $d.Content.Find.Execute("16+2=18", $true, $false, $false, $false, $false, $true, 1, $false, "76-72=4", 2)
$null = $d.Content.Find.Execute("18+39=57", $true, $false, $false, $false, $false, $true, 1, $false, "34+42=76", 2)
$null = $d.Content.Find.Execute("80-49=31", $true, $false, $false, $false, $false, $true, 1, $false, "85-8=77", 2)
$null = $d.Content.Find.Execute("80-75=5", $true, $false, $false, $false, $false, $true, 1, $false, "41-28=13", 2)
$null = $d.Content.Find.Execute("90-64=26", $true, $false, $false, $false, $false, $true, 1, $false, "39+59=98", 2)
$null = $d.Content.Find.Execute("47-42=5", $true, $false, $false, $false, $false, $true, 1, $false, "98-26=72", 2)
$null = $d.Content.Find.Execute("84-38=46", $true, $false, $false, $false, $false, $true, 1, $false, "39+36=75", 2)
$null = $d.Content.Find.Execute("80-48=32", $true, $false, $false, $false, $false, $true, 1, $false, "95-58=37", 2)
$null = $d.Content.Find.Execute("62+5=67", $true, $false, $false, $false, $false, $true, 1, $false, "69-4=65", 2)
$null = $d.Content.Find.Execute("51-48=3", $true, $false, $false, $false, $false, $true, 1, $false, "6+0=6", 2)
$null = $d.Content.Find.Execute("84-59=25", $true, $false, $false, $false, $false, $true, 1, $false, "24+27=51", 2)
$null = $d.Content.Find.Execute("6+60=66", $true, $false, $false, $false, $false, $true, 1, $false, "64-40=24", 2)
$null = $d.Content.Find.Execute("49-48=1", $true, $false, $false, $false, $false, $true, 1, $false, "57-43=14", 2)
$null = $d.Content.Find.Execute("8+78=86", $true, $false, $false, $false, $false, $true, 1, $false, "29+12=41", 2)
$null = $d.Content.Find.Execute("16+32=48", $true, $false, $false, $false, $false, $true, 1, $false, "81+15=96", 2)
$null = $d.Content.Find.Execute("82-13=69", $true, $false, $false, $false, $false, $true, 1, $false, "54-30=24", 2)
$null = $d.Content.Find.Execute("96-61=35", $true, $false, $false, $false, $false, $true, 1, $false, "19+75=94", 2)
$null = $d.Content.Find.Execute("26+44=70", $true, $false, $false, $false, $false, $true, 1, $false, "1+46=47", 2)
$null = $d.Content.Find.Execute("79-46=33", $true, $false, $false, $false, $false, $true, 1, $false, "5+16=21", 2)
$null = $d.Content.Find.Execute("71-28=43", $true, $false, $false, $false, $false, $true, 1, $false, "84-3=81", 2)
$null = $d.Content.Find.Execute("28-6=22", $true, $false, $false, $false, $false, $true, 1, $false, "77+13=90", 2)
$null = $d.Content.Find.Execute("8+66=74", $true, $false, $false, $false, $false, $true, 1, $false, "88-34=54", 2)
$null = $d.Content.Find.Execute("98-31=67", $true, $false, $false, $false, $false, $true, 1, $false, "30-2=28", 2)
$null = $d.Content.Find.Execute("90-25=65", $true, $false, $false, $false, $false, $true, 1, $false, "17-12=5", 2)
$null = $d.Content.Find.Execute("62-30=32", $true, $false, $false, $false, $false, $true, 1, $false, "96-35=61", 2)
$null = $d.Content.Find.Execute("47-43=4", $true, $false, $false, $false, $false, $true, 1, $false, "13+62=75", 2)
$null = $d.Content.Find.Execute("12+76=88", $true, $false, $false, $false, $false, $true, 1, $false, "21+65=86", 2)
$null = $d.Content.Find.Execute("9+64=73", $true, $false, $false, $false, $false, $true, 1, $false, "80+5=85", 2)
$null = $d.Content.Find.Execute("90-88=2", $true, $false, $false, $false, $false, $true, 1, $false, "84+9=93", 2)
$null = $d.Content.Find.Execute("44-30=14", $true, $false, $false, $false, $false, $true, 1, $false, "30-24=6", 2)
$null = $d.Content.Find.Execute("87+12=99", $true, $false, $false, $false, $false, $true, 1, $false, "19+1=20", 2)
$null = $d.Content.Find.Execute("68-14=54", $true, $false, $false, $false, $false, $true, 1, $false, "92-85=7", 2)
$null = $d.Content.Find.Execute("33-26=7", $true, $false, $false, $false, $false, $true, 1, $false, "99-71=28", 2)
$null = $d.Content.Find.Execute("83-35=48", $true, $false, $false, $false, $false, $true, 1, $false, "53-51=2", 2)
$null = $d.Content.Find.Execute("19-4=15", $true, $false, $false, $false, $false, $true, 1, $false, "76-47=29", 2)
$null = $d.Content.Find.Execute("75-11=64", $true, $false, $false, $false, $false, $true, 1, $false, "58-49=9", 2)
$null = $d.Content.Find.Execute("42+13=55", $true, $false, $false, $false, $false, $true, 1, $false, "50+25=75", 2)
$null = $d.Content.Find.Execute("94-80=14", $true, $false, $false, $false, $false, $true, 1, $false, "22+76=98", 2)
$null = $d.Content.Find.Execute("82-43=39", $true, $false, $false, $false, $false, $true, 1, $false, "19+57=76", 2)
$null = $d.Content.Find.Execute("69-30=39", $true, $false, $false, $false, $false, $true, 1, $false, "49-25=24", 2)
$null = $d.Content.Find.Execute("93+1=94", $true, $false, $false, $false, $false, $true, 1, $false, "59-40=19", 2)
